# Completar módulo y lógica cargue proveedores
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$headers = @("NOMBRE", "NIT", "CORREO", "ASESOR CONTACTO", "TELÉFONO", "PRODUCTO O SERVICIO QUE OFRECEN")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 1
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.Font.Color = 16777215
    $cell.Interior.Color = 12611584
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4108
}

# --- Example data row ---
$ws.Cells.Item(2, 1).Value = "7 CUEROS "
$ws.Cells.Item(2, 2).Value = "830101585-1"
$ws.Cells.Item(2, 3).Value = "cilia.coronado@7-cueros.com"
$ws.Cells.Item(2, 4).Value = "CILIA CORONADO "
$ws.Cells.Item(2, 5).Value = 3222274161
$ws.Cells.Item(2, 6).Value = "CALZADO "

$dataRange = $ws.Range("A2:F2")
$dataRange.HorizontalAlignment = -4108
$dataRange.VerticalAlignment = -4108

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 28
$ws.Columns.Item(2).ColumnWidth = 16.88671875
$ws.Columns.Item(3).ColumnWidth = 17.6640625
$ws.Columns.Item(4).ColumnWidth = 19.5546875
$ws.Columns.Item(5).ColumnWidth = 23.88671875
$ws.Columns.Item(6).ColumnWidth = 46.77734375

# --- Workbook metadata ---
$wb.Title = "ejemplo_proveedores"
